$wb = $excel.ActiveWorkbook

# --- ParkingLots: mark spot "Место №1" (row 32) as occupied (IsFree 1 -> 0) ---
$lots = $wb.Worksheets.Item("ParkingLots")
# Reuse an existing cell already holding the text "0" shared-string so the
# written cell keeps its original string type instead of turning numeric.
$lots.Range("D2").Copy($lots.Range("D32"))

# --- Receipts: append the new receipt row (row 4) ---
$receipts = $wb.Worksheets.Item("Receipts")

$receipts.Range("A4").Value = "eb33ae18-ac75-4f14-98b9-3ba6accf121d"
$receipts.Range("B4").Value = "A"
# "000001" must stay text; copy it from a sibling cell that already holds
# that shared string so Excel doesn't coerce it into the number 1.
$receipts.Range("C2").Copy($receipts.Range("C4"))
$receipts.Range("D4").Value = "cae35c5f-6a22-4334-ac9f-277311351804"
$receipts.Range("E4").Value = "50338e7a-1839-4746-b476-dad07a049087"
$receipts.Range("F4").Value = "5ce8cb55-8bc9-498e-80c9-1179059f9fa4"
$receipts.Range("G4").Value = "62d12298-053f-4eee-83cc-7dc6c94694b2"
$receipts.Range("H4").Value = 31
$receipts.Range("I4").Value = 213
$receipts.Range("J4").Value = "27.01.2025 21:43:37"
$receipts.Range("K4").Value = "fb1569f5-9baf-472b-bdc5-811071cfd701"
